$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("T3").Value = 1.04

# Row 4
$ws.Range("I4").Value = 2.46
$ws.Range("Q4").Value = 1.83
$ws.Range("R4").Value = 1.33
$ws.Range("T4").Value = 1.55
$ws.Range("U4").Value = 1.84
$ws.Range("V4").Value = 1.68
$ws.Range("AJ4").Value = 85
$ws.Range("AL4").Value = 65

# Row 5
$ws.Range("H5").Value = 3.55
$ws.Range("J5").Value = 3.25

# Row 6
$ws.Range("T6").Value = 1.04
$ws.Range("V6").Value = 2.12

# Row 7
$ws.Range("F7").Value = 2.14
$ws.Range("I7").Value = 3.95
$ws.Range("R7").Value = 1.24
$ws.Range("T7").Value = 1.04
$ws.Range("V7").Value = 1.35

# Row 8
$ws.Range("H8").Value = 1.09
$ws.Range("I8").Value = 18
$ws.Range("K8").Value = 950
$ws.Range("P8").Value = 2.84
$ws.Range("V8").Value = 1.06
$ws.Range("W8").Value = 3.95

# Row 9
$ws.Range("P9").Value = 1.73
$ws.Range("T9").Value = 2
$ws.Range("Y9").Value = 12

# Row 10
$ws.Range("F10").Value = 5.7
$ws.Range("L10").Value = 1.35
$ws.Range("Q10").Value = 1.86
$ws.Range("Z10").Value = 9.800000000000001
$ws.Range("AC10").Value = 9.4

# Row 11
$ws.Range("L11").Value = 1.27
$ws.Range("U11").Value = 2.42

# Row 12
$ws.Range("P12").Value = 2.18

# Row 14
$ws.Range("F14").Value = 2.44
$ws.Range("G14").Value = 2.82
$ws.Range("R14").Value = 1.29
$ws.Range("T14").Value = 1.78
$ws.Range("V14").Value = 1.44
$ws.Range("X14").Value = 15.5
$ws.Range("Y14").Value = 13.5
$ws.Range("Z14").Value = 25
$ws.Range("AA14").Value = 65
$ws.Range("AB14").Value = 12.5
$ws.Range("AC14").Value = 9.199999999999999
$ws.Range("AD14").Value = 16.5
$ws.Range("AE14").Value = 46
$ws.Range("AF14").Value = 21
$ws.Range("AG14").Value = 15
$ws.Range("AH14").Value = 22
$ws.Range("AI14").Value = 60
$ws.Range("AJ14").Value = 48
$ws.Range("AK14").Value = 38
$ws.Range("AN14").Value = 34
$ws.Range("AO14").Value = 44

# Row 15
$ws.Range("P15").Value = 1.73

# Row 16
$ws.Range("T16").Value = 1.63
$ws.Range("U16").Value = 1.65

# Row 19
$ws.Range("G19").Value = 1.53
$ws.Range("H19").Value = 6.6
$ws.Range("I19").Value = 8.199999999999999
$ws.Range("Q19").Value = 1.47
$ws.Range("R19").Value = 1.63
$ws.Range("S19").Value = 2
$ws.Range("V19").Value = 1.13
$ws.Range("Y19").Value = 50
$ws.Range("AN19").Value = 7

# Row 20
$ws.Range("F20").Value = 1.45
$ws.Range("G20").Value = 1.62
$ws.Range("H20").Value = 5.5
$ws.Range("I20").Value = 8
$ws.Range("J20").Value = 4.1
$ws.Range("Q20").Value = 1.38
$ws.Range("R20").Value = 1.38
$ws.Range("V20").Value = 1.14
$ws.Range("W20").Value = 2.58
$ws.Range("Z20").Value = 85
$ws.Range("AC20").Value = 17
$ws.Range("AG20").Value = 14.5

# Row 22
$ws.Range("I22").Value = 2.7
$ws.Range("J22").Value = 3.55
$ws.Range("P22").Value = 1.93
$ws.Range("Q22").Value = 1.87
$ws.Range("V22").Value = 1.58

# Row 23
$ws.Range("AC23").Value = 13
$ws.Range("AF23").Value = 15.5
$ws.Range("AG23").Value = 15
$ws.Range("AI23").Value = 100
$ws.Range("AJ23").Value = 980
$ws.Range("AN23").Value = 11

# Row 24
$ws.Range("F24").Value = 4.5
$ws.Range("G24").Value = 8.4
$ws.Range("H24").Value = 1.44
$ws.Range("I24").Value = 1.77
$ws.Range("J24").Value = 4.3
$ws.Range("N24").Value = 2.42
$ws.Range("Q24").Value = 1.39

# Row 25
$ws.Range("G25").Value = 2.98
$ws.Range("H25").Value = 2.84
$ws.Range("I25").Value = 3.45
$ws.Range("N25").Value = 1.1
$ws.Range("P25").Value = 1.58
$ws.Range("R25").Value = 1.08
$ws.Range("W25").Value = 1.51

# Row 26
$ws.Range("F26").Value = 1.18
$ws.Range("G26").Value = 1.72
$ws.Range("H26").Value = 5.1
$ws.Range("I26").Value = 9.6
$ws.Range("K26").Value = 1000
$ws.Range("N26").Value = 2.18
$ws.Range("P26").Value = 2.18
$ws.Range("T26").Value = 1.04
$ws.Range("U26").Value = 1.04
$ws.Range("W26").Value = 2.38
$ws.Range("X26").Value = 990
$ws.Range("Y26").Value = 990
$ws.Range("Z26").Value = 1000
$ws.Range("AB26").Value = 990
$ws.Range("AC26").Value = 990
$ws.Range("AD26").Value = 990
$ws.Range("AF26").Value = 980
$ws.Range("AG26").Value = 990
$ws.Range("AH26").Value = 990
$ws.Range("AI26").Value = 1000
$ws.Range("AJ26").Value = 980
$ws.Range("AK26").Value = 980
$ws.Range("AL26").Value = 980
$ws.Range("AN26").Value = 1000

# Row 27
$ws.Range("F27").Value = 2.2
$ws.Range("H27").Value = 2.86
$ws.Range("I27").Value = 3.35
$ws.Range("J27").Value = 3.75
$ws.Range("L27").Value = 1.2
$ws.Range("M27").Value = 1.04
$ws.Range("N27").Value = 5.1
$ws.Range("T27").Value = 1.54
$ws.Range("V27").Value = 1.43
$ws.Range("Z27").Value = 28
$ws.Range("AB27").Value = 15.5
$ws.Range("AD27").Value = 16
$ws.Range("AE27").Value = 980
$ws.Range("AI27").Value = 36
$ws.Range("AO27").Value = 21

# Row 28
$ws.Range("P28").Value = 1.77
$ws.Range("Q28").Value = 1.92
$ws.Range("S28").Value = 2.5
$ws.Range("T28").Value = 1.64
$ws.Range("U28").Value = 1.69
$ws.Range("X28").Value = 17.5
$ws.Range("Y28").Value = 21
$ws.Range("Z28").Value = 48
$ws.Range("AC28").Value = 11.5
$ws.Range("AD28").Value = 26
$ws.Range("AE28").Value = 85
$ws.Range("AF28").Value = 18
$ws.Range("AG28").Value = 15.5
$ws.Range("AI28").Value = 95
$ws.Range("AJ28").Value = 34
$ws.Range("AK28").Value = 32
$ws.Range("AN28").Value = 23

# Row 30
$ws.Range("P30").Value = 2.24
$ws.Range("R30").Value = 1.5

# Row 31
$ws.Range("I31").Value = 8.6
$ws.Range("Q31").Value = 1.94
$ws.Range("V31").Value = 1.13
$ws.Range("AA31").Value = 310
$ws.Range("AG31").Value = 10.5

# Row 32
$ws.Range("O32").Value = 1.2
$ws.Range("R32").Value = 1.61
$ws.Range("AJ32").Value = 11
